$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.431.87"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
$ws.Range("D3").Value = "1.583.81"
$ws.Range("E3").Value = "  -0.21%  "

# Row 6
$ws.Range("E6").Value = "  -0.04%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("D8").Value = "'44.42"
$ws.Range("E8").Value = "  -0.21%  "

# Row 9
$ws.Range("D9").Value = "'24.02"
$ws.Range("E9").Value = "  -0.85%  "

# Row 10
$ws.Range("E10").Value = "  -1.95%  "

# Row 11
$ws.Range("D11").Value = "'0.0591"
$ws.Range("E11").Value = "  -1.46%  "

# Row 12
$ws.Range("E12").Value = "  +1.01%  "

# Row 13
$ws.Range("D13").Value = "1.810.76"

# Row 14
$ws.Range("D14").Value = "1.583.44"
$ws.Range("E14").Value = "  -0.34%  "

# Row 15
$ws.Range("D15").Value = "'3.70"
$ws.Range("E15").Value = "  -0.69%  "

# Row 16
$ws.Range("E16").Value = "  -1.52%  "

# Row 17
$ws.Range("D17").Value = "28.465.79"
$ws.Range("E17").Value = "  -0.30%  "

# Row 18
$ws.Range("D18").Value = "'62.18"
$ws.Range("E18").Value = "  -1.24%  "

# Row 19
$ws.Range("D19").Value = "'230.50"
$ws.Range("E19").Value = "  -0.56%  "

# Row 20
$ws.Range("E20").Value = "  -0.59%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0690"
$ws.Range("E21").Value = "  -2.24%  "

# Row 22
$ws.Range("E22").Value = "  +0.14%  "

# Row 23
$ws.Range("E23").Value = "  -2.96%  "

# Row 24
$ws.Range("D24").Value = "'9.14"
$ws.Range("E24").Value = "  -1.77%  "

# Row 25
$ws.Range("E25").Value = "  +3.74%  "

# Row 26
$ws.Range("D26").Value = "'152.08"
$ws.Range("E26").Value = "  +0.33%  "

# Row 27
$ws.Range("D27").Value = "'15.01"
$ws.Range("E27").Value = "  -1.38%  "

# Row 28
$ws.Range("E28").Value = "  -1.51%  "

# Row 29
$ws.Range("E29").Value = "  -1.76%  "

# Row 30
$ws.Range("E30").Value = "  +0.10%  "

# Row 31
$ws.Range("E31").Value = "  +2.84%  "

# Row 32
$ws.Range("E32").Value = "  -1.67%  "

# Row 33
$ws.Range("E33").Value = "  -1.01%  "

# Row 34
$ws.Range("E34").Value = "  -2.07%  "

# Row 35
$ws.Range("D35").Value = "1.398.39"
$ws.Range("E35").Value = "  +0.76%  "

# Row 36
$ws.Range("D36").Value = "'1.07"

# Row 37
$ws.Range("E37").Value = "  -4.11%  "

# Row 38
$ws.Range("D38").Value = "'2.36"
$ws.Range("E38").Value = "  +0.49%  "

# Row 39
$ws.Range("D39").Value = "'2.66"
$ws.Range("E39").Value = "  +1.59%  "

# Row 40
$ws.Range("E40").Value = "  -0.65%  "

# Row 41
$ws.Range("D41").Value = "'0.523"
$ws.Range("E41").Value = "  -3.23%  "

# Row 42
$ws.Range("E42").Value = "  +0.14%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.791"
$ws.Range("E43").Value = "  -2.46%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.89"
$ws.Range("E44").Value = "  +1.56%  "

# Row 45
$ws.Range("D45").Value = "'0.0462"
$ws.Range("E45").Value = "  -0.71%  "

# Row 46
$ws.Range("D46").Value = "'5.43"
$ws.Range("E46").Value = "  -3.43%  "

# Row 47
$ws.Range("E47").Value = "  -1.99%  "

# Row 48
$ws.Range("D48").Value = "'63.05"
$ws.Range("E48").Value = "  +0.39%  "

# Row 49
$ws.Range("D49").Value = "1.722.06"
$ws.Range("E49").Value = "  -0.16%  "

# Row 50
$ws.Range("D50").Value = "'86.66"
$ws.Range("E50").Value = "  -0.45%  "

# Row 51
$ws.Range("E51").Value = "  -2.16%  "
